$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text entered after the meeting / new experiment, in the order typed
# (this governs the shared-string table ordering)
$ws.Range("D39").Value = "only add energy of node when jobs are running on it."

$ws.Range("L41").Value = "after a certain node and time, not worthwhile to add any resources"
$ws.Range("I41").Value = "time vs energy 10k/20k/30k"
$ws.Range("I40").Value = "graph"
$ws.Range("L40").Value = "message"

$ws.Range("I42").Value = "what are nodes doing?"
$ws.Range("L42").Value = "turn them on and off if not doing anything"

$ws.Range("I44").Value = "number of jobs"

$ws.Range("H51").Value = "outcome: "
$ws.Range("I51").Value = "better scheduler"
$ws.Range("I52").Value = "room for improvement"

$ws.Range("K41").Value = "above has no linear speedup"
$ws.Range("K42").Value = "workflow does not use resources equally"
$ws.Range("R42").Value = "<- nodes are not fully being utilised and hence can be turned off."
$ws.Range("R44").Value = "<- to max out  number of nodes have"
$ws.Range("K43").Value = "two lines of same graph - one with energy consumption of all nodes and second with just the nodes that are working."

$ws.Range("X42").Value = "<- dotted graph showing the nodes are ideal"
$ws.Range("X43").Value = "<- shows that if the nodes are turned off, it will be good"

# Numeric cells in the new "outcome" / step table
$ws.Range("H41").Value = 1
$ws.Range("H42").Value = 2.1
$ws.Range("H43").Value = 2.2
$ws.Range("H44").Value = 3

# Update the view state to match where the author ended up after the edits
$ws.Range("N39").Select()
